$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44301, 1, 10, 99.30486593843098),
    @(44302, 3, 12, 119.1658391261172),
    @(44303, 4, 13, 129.0963257199603)
)

$startRow = 227
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$ws.Range("A226").Copy() | Out-Null
$ws.Range("A227:A229").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
